$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: Acierto, profit = cuota - 1 = 1.38
$ws.Range("G22").Value = "Acierto"
$ws.Range("H22").Value = 1.38

# Row 38: Fallo, profit = -1
$ws.Range("G38").Value = "Fallo"
$ws.Range("H38").Value = -1

# Row 39: Fallo, profit = -1
$ws.Range("G39").Value = "Fallo"
$ws.Range("H39").Value = -1

# Row 42: Fallo, profit = -1
$ws.Range("G42").Value = "Fallo"
$ws.Range("H42").Value = -1

# Row 44: Acierto, profit = 1.1
$ws.Range("G44").Value = "Acierto"
$ws.Range("H44").Value = 1.1

# Row 50: Fallo, profit = -1
$ws.Range("G50").Value = "Fallo"
$ws.Range("H50").Value = -1

# Row 63: Acierto, profit = 0.53
$ws.Range("G63").Value = "Acierto"
$ws.Range("H63").Value = 0.53

# Row 64: Fallo, profit = -1
$ws.Range("G64").Value = "Fallo"
$ws.Range("H64").Value = -1

# Row 69: A69 event_id becomes a true number instead of text
$ws.Range("A69").Value = 14601390

# Row 70: A70 event_id becomes a true number instead of text
$ws.Range("A70").Value = 14601341
